$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") is bumped by one day (46062 -> 46063) for every data row.
foreach ($r in 2..13) {
    $ws.Cells.Item($r, 3).Value = 46063
}

# Rows 4-10 get reshuffled (same 7 records, different row order).
# Capture the "before" values for A, B, G and the hyperlink formulas (S,T,V,W,X,Y)
# for each of those rows before overwriting anything.
$rows = 4..10
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        A = $ws.Cells.Item($r, 1).Value()
        B = $ws.Cells.Item($r, 2).Value()
        G = $ws.Cells.Item($r, 7).Value()
        S = $ws.Cells.Item($r, 19).Formula()
        T = $ws.Cells.Item($r, 20).Formula()
        V = $ws.Cells.Item($r, 22).Formula()
        W = $ws.Cells.Item($r, 23).Formula()
        X = $ws.Cells.Item($r, 24).Formula()
        Y = $ws.Cells.Item($r, 25).Formula()
    }
}

# New row order: destination row -> source (old) row
$order = @{ 4 = 5; 5 = 4; 6 = 10; 7 = 9; 8 = 7; 9 = 8; 10 = 6 }

foreach ($dest in $rows) {
    $src = $order[$dest]
    $rec = $data[$src]

    $ws.Cells.Item($dest, 1).Value = $rec.A
    $ws.Cells.Item($dest, 2).Value = $rec.B
    $ws.Cells.Item($dest, 7).Value = $rec.G

    if ($rec.S) { $ws.Cells.Item($dest, 19).Formula = $rec.S }
    if ($rec.T) { $ws.Cells.Item($dest, 20).Formula = $rec.T }
    if ($rec.V) { $ws.Cells.Item($dest, 22).Formula = $rec.V }
    if ($rec.W) { $ws.Cells.Item($dest, 23).Formula = $rec.W }
    if ($rec.X) { $ws.Cells.Item($dest, 24).Formula = $rec.X }
    if ($rec.Y) { $ws.Cells.Item($dest, 25).Formula = $rec.Y }
}
